# Natmi following Dr Hou advice
# Update recomputed NATMI ligand/receptor-expressing cell counts and all
# derived expression / specificity / edge-weight metrics for the
# Lrpap1-Lrp1 sheet (Ligand/Receptor-expressing cells now 3 instead of 1,
# with corresponding totals, specificities and edge weights recomputed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.660188000000001
$ws.Range("H2").Value = 19.980564
$ws.Range("I2").Value = 0.1500148400131262
$ws.Range("J2").Value = 0.1500148400131261
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 157.4188559899613
$ws.Range("R2").Value = 1416.769703909652
$ws.Range("S2").Value = 0.0102427801344648
$ws.Range("T2").Value = 0.0102427801344648

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.660188000000001
$ws.Range("H3").Value = 19.980564
$ws.Range("I3").Value = 0.1500148400131262
$ws.Range("J3").Value = 0.1500148400131261
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 1207.414768497571
$ws.Range("R3").Value = 10866.73291647814
$ws.Range("S3").Value = 0.07856291374404983
$ws.Range("T3").Value = 0.07856291374404981

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.660188000000001
$ws.Range("H4").Value = 19.980564
$ws.Range("I4").Value = 0.1500148400131262
$ws.Range("J4").Value = 0.1500148400131261
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 739.9505254827106
$ws.Range("R4").Value = 6659.554729344396
$ws.Range("S4").Value = 0.04814639577475029
$ws.Range("T4").Value = 0.04814639577475028

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.660188000000001
$ws.Range("H5").Value = 19.980564
$ws.Range("I5").Value = 0.1500148400131262
$ws.Range("J5").Value = 0.1500148400131261
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 200.7583088513947
$ws.Range("R5").Value = 1806.824779662552
$ws.Range("S5").Value = 0.01306275035986125
$ws.Range("T5").Value = 0.01306275035986125

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.410331
$ws.Range("H6").Value = 49.230993
$ws.Range("I6").Value = 0.3696281815959916
$ws.Range("J6").Value = 0.3696281815959916
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 387.8712631590276
$ws.Range("R6").Value = 3490.841368431249
$ws.Range("S6").Value = 0.02523763779142448
$ws.Range("T6").Value = 0.02523763779142448

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.410331
$ws.Range("H7").Value = 49.230993
$ws.Range("I7").Value = 0.3696281815959916
$ws.Range("J7").Value = 0.3696281815959916
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 2975.002508237531
$ws.Range("R7").Value = 26775.02257413778
$ws.Range("S7").Value = 0.1935746286537718
$ws.Range("T7").Value = 0.1935746286537717

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 16.410331
$ws.Range("H8").Value = 49.230993
$ws.Range("I8").Value = 0.3696281815959916
$ws.Range("J8").Value = 0.3696281815959916
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 1823.196739610836
$ws.Range("R8").Value = 16408.77065649752
$ws.Range("S8").Value = 0.1186300283296288
$ws.Range("T8").Value = 0.1186300283296288

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 16.410331
$ws.Range("H9").Value = 49.230993
$ws.Range("I9").Value = 0.3696281815959916
$ws.Range("J9").Value = 0.3696281815959916
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 494.6572528060192
$ws.Range("R9").Value = 4451.915275254173
$ws.Range("S9").Value = 0.03218588682116665
$ws.Range("T9").Value = 0.03218588682116665

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.478895
$ws.Range("H10").Value = 43.436685
$ws.Range("I10").Value = 0.3261242951387937
$ws.Range("J10").Value = 0.3261242951387937
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 342.2202326569117
$ws.Range("R10").Value = 3079.982093912205
$ws.Range("S10").Value = 0.02226726003455183
$ws.Range("T10").Value = 0.02226726003455183

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 14.478895
$ws.Range("H11").Value = 43.436685
$ws.Range("I11").Value = 0.3261242951387937
$ws.Range("J11").Value = 0.3261242951387937
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 2624.855582834244
$ws.Range("R11").Value = 23623.70024550819
$ws.Range("S11").Value = 0.1707916021280712
$ws.Range("T11").Value = 0.1707916021280712

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 14.478895
$ws.Range("H12").Value = 43.436685
$ws.Range("I12").Value = 0.3261242951387937
$ws.Range("J12").Value = 0.3261242951387937
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 1608.613144802968
$ws.Range("R12").Value = 14477.51830322671
$ws.Range("S12").Value = 0.104667707435744
$ws.Range("T12").Value = 0.104667707435744

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 14.478895
$ws.Range("H13").Value = 43.436685
$ws.Range("I13").Value = 0.3261242951387937
$ws.Range("J13").Value = 0.3261242951387937
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 436.4379014882033
$ws.Range("R13").Value = 3927.941113393829
$ws.Range("S13").Value = 0.02839772554042667
$ws.Range("T13").Value = 0.02839772554042667

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.847447
$ws.Range("H14").Value = 20.542341
$ws.Range("I14").Value = 0.1542326832520885
$ws.Range("J14").Value = 0.1542326832520885
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 161.8448718252237
$ws.Range("R14").Value = 1456.603846427013
$ws.Range("S14").Value = 0.01053076791577063
$ws.Range("T14").Value = 0.01053076791577063

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.847447
$ws.Range("H15").Value = 20.542341
$ws.Range("I15").Value = 0.1542326832520885
$ws.Range("J15").Value = 0.1542326832520885
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 1241.36265137026
$ws.Range("R15").Value = 11172.26386233233
$ws.Range("S15").Value = 0.0807718022416113
$ws.Range("T15").Value = 0.08077180224161129

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.847447
$ws.Range("H16").Value = 20.542341
$ws.Range("I16").Value = 0.1542326832520885
$ws.Range("J16").Value = 0.1542326832520885
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 760.7551026885443
$ws.Range("R16").Value = 6846.795924196898
$ws.Range("S16").Value = 0.04950008818198923
$ws.Range("T16").Value = 0.04950008818198923

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.847447
$ws.Range("H17").Value = 20.542341
$ws.Range("I17").Value = 0.1542326832520885
$ws.Range("J17").Value = 0.1542326832520885
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 206.4028642539153
$ws.Range("R17").Value = 1857.625778285238
$ws.Range("S17").Value = 0.0134300249127173
$ws.Range("T17").Value = 0.0134300249127173
